$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets contain the same duplicated data rows (F2, F3, F5)
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 707
    $ws.Range("F3").Value = 4048
    $ws.Range("F5").Value = 744
}
